$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values in columns A, B, E, F, G, H, Q, R between row 2 and row 3.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"

    $cell2 = $ws.Range($addr2)
    $cell3 = $ws.Range($addr3)

    $val2 = $cell2.Value2
    $val3 = $cell3.Value2

    $cell2.Value2 = $val3
    $cell3.Value2 = $val2
}
